$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, shifting existing rows 193..278 down to 194..279
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new data record
$ws.Cells.Item(193, 1).Value = 11
$ws.Cells.Item(193, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(193, 3).Value = "Bíobío"
$ws.Cells.Item(193, 4).Value = 44489
$ws.Cells.Item(193, 5).Value = 8
$ws.Cells.Item(193, 6).Value = "Fruta"
$ws.Cells.Item(193, 7).Value = 100102
$ws.Cells.Item(193, 8).Value = "Cítricos"
$ws.Cells.Item(193, 9).Value = 100102003
$ws.Cells.Item(193, 10).Value = "Limón"
$ws.Cells.Item(193, 11).Value = "Sin especificar"
$ws.Cells.Item(193, 12).Value = "1a amarillo"
$ws.Cells.Item(193, 13).Value = 650
$ws.Cells.Item(193, 14).Value = 4500
$ws.Cells.Item(193, 15).Value = 5000
$ws.Cells.Item(193, 16).Value = 4769
$ws.Cells.Item(193, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(193, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(193, 19).Value = 265
$ws.Cells.Item(193, 20).Value = 18
